# Versi 1.0.0.22 add ctkHasil.rpt add ctkHasil.vb Cetak Hasil diganti ke CR.
# Populate two extra copies of the FPP (sheet2) data table into columns F:I and K:N
# so the new Crystal-Report binding ("ctkHasil.rpt") has duplicate/aux data regions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FPP")

# Row data: row -> (FG label/test, H qty, I value, KL label/test, M qty, N value)
# FG/H/I block mirrors A:D labels with its own qty/value figures; K:N block is an
# exact duplicate of the original A:D figures.
$rows = @(
    @{ r=9;  F='Air Limbah'; G='BILIRUBIN_KONVENSIONAL';       H=2; I=15000;  K='Air Limbah'; L='BILIRUBIN_KONVENSIONAL';       M=2; N=15000  },
    @{ r=10; F='Air Minum';  G='ALUMINIUM_AIR MINUM_SPEKTRO';  H=9; I=405000; K='Air Minum';  L='ALUMINIUM_AIR MINUM_SPEKTRO';  M=3; N=135000 },
    @{ r=11; F='Air Minum';  G='BILIRUBIN_KONVENSIONAL';       H=4; I=30000;  K='Air Minum';  L='BILIRUBIN_KONVENSIONAL';       M=4; N=30000  },
    @{ r=12; F='AMDK';       G='ALUMINIUM_AIR MINUM_SPEKTRO';  H=9; I=405000; K='AMDK';       L='ALUMINIUM_AIR MINUM_SPEKTRO';  M=3; N=135000 },
    @{ r=13; F='AMDK';       G='BILIRUBIN_KONVENSIONAL';       H=4; I=30000;  K='AMDK';       L='BILIRUBIN_KONVENSIONAL';       M=4; N=30000  },
    @{ r=14; F='apa aja';    G='ALUMINIUM_AIR MINUM_SPEKTRO';  H=6; I=270000; K='apa aja';    L='ALUMINIUM_AIR MINUM_SPEKTRO';  M=2; N=90000  },
    @{ r=15; F='apa aja';    G='BILIRUBIN_KONVENSIONAL';       H=2; I=15000;  K='apa aja';    L='BILIRUBIN_KONVENSIONAL';       M=2; N=15000  },
    @{ r=16; F='okok';       G='ALUMINIUM_AIR MINUM_SPEKTRO';  H=3; I=135000; K='okok';       L='ALUMINIUM_AIR MINUM_SPEKTRO';  M=1; N=45000  },
    @{ r=17; F='okok';       G='BILIRUBIN_KONVENSIONAL';       H=1; I=7500;   K='okok';       L='BILIRUBIN_KONVENSIONAL';       M=1; N=7500   },
    @{ r=18; F='Sperma';     G='ALUMINIUM_AIR MINUM_SPEKTRO';  H=3; I=135000; K='Sperma';     L='ALUMINIUM_AIR MINUM_SPEKTRO';  M=1; N=45000  },
    @{ r=19; F='Sperma';     G='BILIRUBIN_KONVENSIONAL';       H=1; I=7500;   K='Sperma';     L='BILIRUBIN_KONVENSIONAL';       M=1; N=7500   },
    @{ r=20; F='Timun';      G='ALUMINIUM_AIR MINUM_SPEKTRO';  H=6; I=270000; K='Timun';      L='ALUMINIUM_AIR MINUM_SPEKTRO';  M=2; N=90000  },
    @{ r=21; F='Timun';      G='BILIRUBIN_KONVENSIONAL';       H=2; I=15000;  K='Timun';      L='BILIRUBIN_KONVENSIONAL';       M=2; N=15000  },
    @{ r=22; F='Tinja';      G='ALUMINIUM_AIR MINUM_SPEKTRO';  H=3; I=135000; K='Tinja';      L='ALUMINIUM_AIR MINUM_SPEKTRO';  M=1; N=45000  },
    @{ r=23; F='Tinja';      G='BILIRUBIN_KONVENSIONAL';       H=3; I=22500;  K='Tinja';      L='BILIRUBIN_KONVENSIONAL';       M=3; N=22500  }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 6).Value  = $row.F    # F
    $ws.Cells.Item($r, 7).Value  = $row.G    # G
    $ws.Cells.Item($r, 8).Value  = $row.H    # H
    $ws.Cells.Item($r, 9).Value  = $row.I    # I
    $ws.Cells.Item($r, 11).Value = $row.K    # K
    $ws.Cells.Item($r, 12).Value = $row.L    # L
    $ws.Cells.Item($r, 13).Value = $row.M    # M
    $ws.Cells.Item($r, 14).Value = $row.N    # N

    # N column mirrors D's Comma number format
    $ws.Range("N$r").NumberFormat = $ws.Range("D$r").NumberFormat
}

# Totals row 24: C24 / I24 plain sums; M24:N24 as a shared formula (mirrors D24)
$ws.Range("C24").Formula = "=SUM(C9:C23)"
$ws.Range("I24").Formula = "=SUM(I9:I23)"
$ws.Range("M24:N24").Formula = "=SUM(M9:M23)"
$ws.Range("N24").NumberFormat = $ws.Range("D24").NumberFormat

# Column widths for the new K/L/N columns. K mirrors A's font (new content), L/N
# reuse the stored widths of B/D (authored widths: 32.42578125 / 11.5703125).
$ws.Columns.Item(11).ColumnWidth = 9.59
$ws.Columns.Item(12).ColumnWidth = 31.59
$ws.Columns.Item(14).ColumnWidth = 10.74

# Restore selection close to the authored state (engine does not support true
# multi-area sqref selections, so land on the final selected cell, M22).
$null = $ws.Range("M22").Select()
